# Auto-generated script: apply scheduled price-runner updates to Pandaemonium_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Cell="H125"; Value=5208.8184},
    @{Cell="I125"; Value=200},
    @{Cell="J125"; Value=5447.3335},
    @{Cell="K125"; Value=1800},
    @{Cell="L125"; Value=49026.0015},
    @{Cell="M125"; Value=660},
    @{Cell="N125"; Value=-53946.0015}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Cell="H61"; Value=8115.1465},
    @{Cell="I61"; Value=4275.6895},
    @{Cell="J61"; Value=17393.834},
    @{Cell="K61"; Value=4275.6895},
    @{Cell="L61"; Value=17393.834},
    @{Cell="M61"; Value=-4063.6895},
    @{Cell="N61"; Value=-17817.834},
    @{Cell="H74"; Value=1658.0541},
    @{Cell="I74"; Value=1763.68},
    @{Cell="J74"; Value=1438},
    @{Cell="K74"; Value=1763.68},
    @{Cell="L74"; Value=1438},
    @{Cell="M74"; Value=-889.6800000000001},
    @{Cell="N74"; Value=-3186},
    @{Cell="H77"; Value=1658.0541},
    @{Cell="I77"; Value=1763.68},
    @{Cell="J77"; Value=1438},
    @{Cell="K77"; Value=8818.4},
    @{Cell="L77"; Value=7190},
    @{Cell="M77"; Value=-4450.4},
    @{Cell="N77"; Value=-15926},
    @{Cell="H122"; Value=1969.9546},
    @{Cell="I122"; Value=1666.1875},
    @{Cell="J122"; Value=2780},
    @{Cell="K122"; Value=4998.5625},
    @{Cell="L122"; Value=8340},
    @{Cell="M122"; Value=-2548.5625},
    @{Cell="N122"; Value=-13240},
    @{Cell="H136"; Value=8115.1465},
    @{Cell="I136"; Value=4275.6895},
    @{Cell="J136"; Value=17393.834},
    @{Cell="K136"; Value=12827.0685},
    @{Cell="L136"; Value=52181.50199999999},
    @{Cell="M136"; Value=-10277.0685},
    @{Cell="N136"; Value=-57281.50199999999}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{Cell="H99"; Value=1400},
    @{Cell="I99"; Value=1400},
    @{Cell="K99"; Value=1400},
    @{Cell="M99"; Value=98},
    @{Cell="H134"; Value=24369.244},
    @{Cell="I134"; Value=2011.4857},
    @{Cell="K134"; Value=6034.4571},
    @{Cell="M134"; Value=-3499.4571}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Cell="H31"; Value=4756.0713},
    @{Cell="I31"; Value=5491.731},
    @{Cell="J31"; Value=3560.625},
    @{Cell="K31"; Value=5491.731},
    @{Cell="L31"; Value=3560.625},
    @{Cell="M31"; Value=-5196.731},
    @{Cell="N31"; Value=-4150.625},
    @{Cell="H34"; Value=4756.0713},
    @{Cell="I34"; Value=5491.731},
    @{Cell="J34"; Value=3560.625},
    @{Cell="K34"; Value=5491.731},
    @{Cell="L34"; Value=3560.625},
    @{Cell="M34"; Value=-5289.731},
    @{Cell="N34"; Value=-3964.625},
    @{Cell="H58"; Value=2333796.8},
    @{Cell="I58"; Value=4547685.5},
    @{Cell="J58"; Value=3387.158},
    @{Cell="K58"; Value=4547685.5},
    @{Cell="L58"; Value=3387.158},
    @{Cell="M58"; Value=-4547482.5},
    @{Cell="N58"; Value=-3793.158},
    @{Cell="H82"; Value=32590.5},
    @{Cell="J82"; Value=32590.5},
    @{Cell="L82"; Value=32590.5},
    @{Cell="N82"; Value=-33312.5},
    @{Cell="H85"; Value=32590.5},
    @{Cell="J85"; Value=32590.5},
    @{Cell="L85"; Value=32590.5},
    @{Cell="N85"; Value=-35086.5},
    @{Cell="H122"; Value=11959.777},
    @{Cell="I122"; Value=7123.8184},
    @{Cell="J122"; Value=19559.143},
    @{Cell="K122"; Value=21371.4552},
    @{Cell="L122"; Value=58677.429},
    @{Cell="M122"; Value=-18921.4552},
    @{Cell="N122"; Value=-63577.429},
    @{Cell="H134"; Value=2552.6333},
    @{Cell="I134"; Value=1453.7435},
    @{Cell="J134"; Value=4593.4287},
    @{Cell="K134"; Value=4361.2305},
    @{Cell="L134"; Value=13780.2861},
    @{Cell="M134"; Value=-1826.2305},
    @{Cell="N134"; Value=-18850.2861},
    @{Cell="H136"; Value=2333796.8},
    @{Cell="I136"; Value=4547685.5},
    @{Cell="J136"; Value=3387.158},
    @{Cell="K136"; Value=13643056.5},
    @{Cell="L136"; Value=10161.474},
    @{Cell="M136"; Value=-13640506.5},
    @{Cell="N136"; Value=-15261.474},
    @{Cell="H141"; Value=26239.6},
    @{Cell="J141"; Value=33032.285},
    @{Cell="L141"; Value=33032.285},
    @{Cell="N141"; Value=-43392.285}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Cell="H86"; Value=733.6667},
    @{Cell="I86"; Value=502},
    @{Cell="J86"; Value=780},
    @{Cell="K86"; Value=1506},
    @{Cell="L86"; Value=2340},
    @{Cell="M86"; Value=-320},
    @{Cell="N86"; Value=-4712},
    @{Cell="H89"; Value=733.6667},
    @{Cell="I89"; Value=502},
    @{Cell="J89"; Value=780},
    @{Cell="K89"; Value=4518},
    @{Cell="L89"; Value=7020},
    @{Cell="M89"; Value=1410},
    @{Cell="N89"; Value=-18876},
    @{Cell="H122"; Value=1120.5294},
    @{Cell="I122"; Value=963.3333},
    @{Cell="J122"; Value=1154.2142},
    @{Cell="K122"; Value=8669.9997},
    @{Cell="L122"; Value=10387.9278},
    @{Cell="M122"; Value=-6219.9997},
    @{Cell="N122"; Value=-15287.9278}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Cell="H80"; Value=6382},
    @{Cell="I80"; Value=10124.667},
    @{Cell="J80"; Value=3575},
    @{Cell="K80"; Value=10124.667},
    @{Cell="L80"; Value=3575},
    @{Cell="M80"; Value=-9126.666999999999},
    @{Cell="N80"; Value=-5571},
    @{Cell="H83"; Value=6382},
    @{Cell="I83"; Value=10124.667},
    @{Cell="J83"; Value=3575},
    @{Cell="K83"; Value=50623.335},
    @{Cell="L83"; Value=17875},
    @{Cell="M83"; Value=-45631.335},
    @{Cell="N83"; Value=-27859},
    @{Cell="H86"; Value=36143},
    @{Cell="J86"; Value=36143},
    @{Cell="L86"; Value=36143},
    @{Cell="N86"; Value=-38515},
    @{Cell="H89"; Value=36143},
    @{Cell="J89"; Value=36143},
    @{Cell="L89"; Value=108429},
    @{Cell="N89"; Value=-120285},
    @{Cell="H122"; Value=14093.786},
    @{Cell="I122"; Value=13975.75},
    @{Cell="K122"; Value=41927.25},
    @{Cell="M122"; Value=-39477.25},
    @{Cell="H132"; Value=4611.7896},
    @{Cell="I132"; Value=1779.4546},
    @{Cell="J132"; Value=23305.2},
    @{Cell="K132"; Value=5338.3638},
    @{Cell="L132"; Value=69915.60000000001},
    @{Cell="M132"; Value=-2808.3638},
    @{Cell="N132"; Value=-74975.60000000001}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Cell="H16"; Value=1866.7273},
    @{Cell="I16"; Value=1400.4},
    @{Cell="J16"; Value=2255.3333},
    @{Cell="K16"; Value=1400.4},
    @{Cell="L16"; Value=2255.3333},
    @{Cell="M16"; Value=-1230.4},
    @{Cell="N16"; Value=-2595.3333},
    @{Cell="H40"; Value=3650.5625},
    @{Cell="I40"; Value=3740.4},
    @{Cell="J40"; Value=3500.8333},
    @{Cell="K40"; Value=3740.4},
    @{Cell="L40"; Value=3500.8333},
    @{Cell="M40"; Value=-3604.4},
    @{Cell="N40"; Value=-3772.8333},
    @{Cell="H42"; Value=29800},
    @{Cell="J42"; Value=29800},
    @{Cell="L42"; Value=29800},
    @{Cell="N42"; Value=-30926},
    @{Cell="H49"; Value=29800},
    @{Cell="J49"; Value=29800},
    @{Cell="L49"; Value=29800},
    @{Cell="N49"; Value=-30094},
    @{Cell="H61"; Value=731035.2},
    @{Cell="I61"; Value=24865.334},
    @{Cell="J61"; Value=2002141},
    @{Cell="K61"; Value=24865.334},
    @{Cell="L61"; Value=2002141},
    @{Cell="M61"; Value=-24663.334},
    @{Cell="N61"; Value=-2002545},
    @{Cell="H93"; Value=0},
    @{Cell="J93"; Value=0},
    @{Cell="L93"; Value=0},
    @{Cell="N93"; Value=$null},
    @{Cell="H113"; Value=731035.2},
    @{Cell="I113"; Value=24865.334},
    @{Cell="J113"; Value=2002141},
    @{Cell="K113"; Value=24865.334},
    @{Cell="L113"; Value=2002141},
    @{Cell="M113"; Value=-22695.334},
    @{Cell="N113"; Value=-2006481},
    @{Cell="H132"; Value=5036.35},
    @{Cell="I132"; Value=5491.1396},
    @{Cell="J132"; Value=3886},
    @{Cell="K132"; Value=16473.4188},
    @{Cell="L132"; Value=11658},
    @{Cell="M132"; Value=-13943.4188},
    @{Cell="N132"; Value=-16718},
    @{Cell="H136"; Value=4152.418},
    @{Cell="I136"; Value=2488.7666},
    @{Cell="J136"; Value=6148.8},
    @{Cell="K136"; Value=7466.2998},
    @{Cell="L136"; Value=18446.4},
    @{Cell="M136"; Value=-4916.2998},
    @{Cell="N136"; Value=-23546.4}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Cell="H82"; Value=39856.223},
    @{Cell="I82"; Value=0},
    @{Cell="J82"; Value=39856.223},
    @{Cell="K82"; Value=0},
    @{Cell="L82"; Value=39856.223},
    @{Cell="M82"; Value=$null},
    @{Cell="N82"; Value=-40622.223},
    @{Cell="H85"; Value=39856.223},
    @{Cell="I85"; Value=0},
    @{Cell="J85"; Value=39856.223},
    @{Cell="K85"; Value=0},
    @{Cell="L85"; Value=39856.223},
    @{Cell="M85"; Value=$null},
    @{Cell="N85"; Value=-42508.223},
    @{Cell="H122"; Value=12100},
    @{Cell="I122"; Value=3250},
    @{Cell="J122"; Value=18000},
    @{Cell="K122"; Value=9750},
    @{Cell="L122"; Value=54000},
    @{Cell="M122"; Value=-7300},
    @{Cell="N122"; Value=-58900},
    @{Cell="H123"; Value=28944},
    @{Cell="J123"; Value=28944},
    @{Cell="L123"; Value=28944},
    @{Cell="N123"; Value=-38744},
    @{Cell="H136"; Value=4236.1924},
    @{Cell="I136"; Value=4226.75},
    @{Cell="J136"; Value=4257.4375},
    @{Cell="K136"; Value=12680.25},
    @{Cell="L136"; Value=12772.3125},
    @{Cell="M136"; Value=-10130.25},
    @{Cell="N136"; Value=-17872.3125}
)
foreach ($u in $updates) {
    if ($null -eq $u.Value) {
        $ws.Range($u.Cell).ClearContents()
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
